$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F ("想去人数") on several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 4116
$ws1.Range("F4").Value  = 2383
$ws1.Range("F5").Value  = 474
$ws1.Range("F7").Value  = 38
$ws1.Range("F8").Value  = 35
$ws1.Range("F11").Value = 100
$ws1.Range("F12").Value = 143
$ws1.Range("F13").Value = 1545
$ws1.Range("F14").Value = 280
$ws1.Range("F15").Value = 3025
$ws1.Range("F16").Value = 207

# Sheet "全部类型" (sheet4): update column F ("想去人数") on several rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 4116
$ws4.Range("F4").Value  = 2383
$ws4.Range("F5").Value  = 474
$ws4.Range("F8").Value  = 38
$ws4.Range("F9").Value  = 35
$ws4.Range("F13").Value = 100
$ws4.Range("F14").Value = 143
$ws4.Range("F17").Value = 1546
$ws4.Range("F18").Value = 280
$ws4.Range("F19").Value = 3025
$ws4.Range("F20").Value = 207

$wb.Save()
